$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Rongfeng%Qi%NULL%0,                    Wei%Chen%NULL%0,                    Saiduo%Liu%NULL%0,                    Paul M.%Thompson%NULL%0,                    Long Jiang%Zhang%NULL%0,                    Fei%Xia%NULL%0,                    Fang%Cheng%NULL%0,                    Ailing%Hong%NULL%0,                    Wesley%Surento%NULL%0,                    Song%Luo%NULL%0,                    Zhi Yuan%Sun%NULL%0,                    Chang Sheng%Zhou%NULL%0,                    Lingjiang%Li%NULL%0,                    Xiangao%Jiang%NULL%0,                    Guang Ming%Lu%NULL%0,                     R.%Qi%null%0,                     W.% Chen%null%0,                     S.% Liu%null%0,                     P. M.% Thompson%null%0,                     L. J.% Zhang%null%0,                     F.% Xia%null%0,                     F.% Cheng%null%0,                     A.% Hong%null%0,                     W.% Surento%null%0,                     S.% Luo%null%0,                     Z. Y.% Sun%null%0,                     C. S.% Zhou%null%0,                     L.% Li%null%0,                     X.% Jiang%null%0,                     G. M. % Lu%null%0,                 R.%Qi%null%1,                 W.% Chen%null%1,                 S.% Liu%null%1,                 P. M.% Thompson%null%1,                 L. J.% Zhang%null%1,                 F.% Xia%null%1,                 F.% Cheng%null%1,                 A.% Hong%null%1,                 W.% Surento%null%1,                 S.% Luo%null%1,                 Z. Y.% Sun%null%1,                 C. S.% Zhou%null%1,                 L.% Li%null%1,                 X.% Jiang%null%1,                 G. M. % Lu%null%1]'
$ws.Range("I2").Value = '_PMC_CROSSREF'

$ws.Range("E3").Value = '[Ömer%Karadaş%NULL%1,                    Bilgin%Öztürk%Drbilgin@gmail.com%2,                    Bilgin%Öztürk%Drbilgin@gmail.com%0,                    Ali Rıza%Sonkaya%NULL%2,                    Ali Rıza%Sonkaya%NULL%0]'
$ws.Range("I3").Value = '_PMC_Springer_CROSSREF'

$ws.Range("E4").Value = '[Claudio%Liguori%NULL%0,                    Mariangela%Pierantozzi%NULL%2,                    Matteo%Spanetta%NULL%2,                    Loredana%Sarmati%NULL%2,                    Novella%Cesta%NULL%2,                    Marco%Iannetta%NULL%2,                    Josuel%Ora%NULL%2,                    Grazia Genga%Mina%NULL%2,                    Ermanno%Puxeddu%NULL%2,                    Ottavia%Balbi%NULL%2,                    Gabriella%Pezzuto%NULL%2,                    Andrea%Magrini%NULL%2,                    Paola%Rogliani%NULL%2,                    Massimo%Andreoni%NULL%2,                    Nicola Biagio%Mercuri%NULL%2]'
$ws.Range("I4").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E5").Value = '[Claudio%Liguori%NULL%0,                    Mariangela%Pierantozzi%NULL%0,                    Matteo%Spanetta%NULL%0,                    Loredana%Sarmati%NULL%0,                    Novella%Cesta%NULL%0,                    Marco%Iannetta%NULL%0,                    Josuel%Ora%NULL%0,                    Grazia Genga%Mina%NULL%0,                    Ermanno%Puxeddu%NULL%0,                    Ottavia%Balbi%NULL%0,                    Gabriella%Pezzuto%NULL%0,                    Andrea%Magrini%NULL%0,                    Paola%Rogliani%NULL%0,                    Massimo%Andreoni%NULL%0,                    Nicola Biagio%Mercuri%NULL%0]'
$ws.Range("I5").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E6").Value = '[Qian%Guo%NULL%1,                    Yuchen%Zheng%NULL%1,                    Jia%Shi%NULL%1,                    Jijun%Wang%NULL%1,                    Guanjun%Li%NULL%1,                    Chunbo%Li%NULL%1,                    John A.%Fromson%NULL%1,                    Yong%Xu%NULL%1,                    Xiaohua%Liu%NULL%1,                    Hua%Xu%NULL%1,                    Tianhong%Zhang%NULL%1,                    Yunfei%Lu%NULL%0,                    Xiaorong%Chen%NULL%1,                    Hao%Hu%NULL%1,                    Yingying%Tang%NULL%1,                    Shuwen%Yang%NULL%1,                    Han%Zhou%NULL%1,                    Xiaoliang%Wang%NULL%1,                    Haiying%Chen%NULL%1,                    Zhen%Wang%NULL%1,                    Zongguo%Yang%NULL%1]'
$ws.Range("I6").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E7").Value = '[ X.%Kong%null%1,                     K.% Zheng%null%1,                     M.% Tang%null%1,                     F.% Kong%null%1,                     J.% Zhou%null%1,                     L.% Diao%null%1,                     S.% Wu%null%1,                     P.% Jiao%null%1,                     T.% Su%null%1,                     Y. % Dong%null%1]'
$ws.Range("I7").Value = '_MedBiorxiv_CROSSREF'

$ws.Range("I8").Value = '_CROSSREF'

$ws.Range("E9").Value = '[Clara%Paz%clara.paz@udla.edu.ec%1,                    Guido%Mascialino%NULL%2,                    Guido%Mascialino%NULL%0,                    Lila%Adana‐Díaz%NULL%2,                    Lila%Adana‐Díaz%NULL%0,                    Alberto%Rodríguez‐Lorenzana%NULL%2,                    Alberto%Rodríguez‐Lorenzana%NULL%0,                    Katherine%Simbaña‐Rivera%NULL%2,                    Katherine%Simbaña‐Rivera%NULL%0,                    Lenin%Gómez‐Barreno%NULL%2,                    Lenin%Gómez‐Barreno%NULL%0,                    Maritza%Troya%NULL%2,                    Maritza%Troya%NULL%0,                    María Ignacia%Páez%NULL%2,                    María Ignacia%Páez%NULL%0,                    Javier%Cárdenas%NULL%2,                    Javier%Cárdenas%NULL%0,                    Rebekka M%Gerstner%NULL%2,                    Rebekka M%Gerstner%NULL%0,                    Esteban%Ortiz‐Prado%NULL%2,                    Esteban%Ortiz‐Prado%NULL%0]'
$ws.Range("I9").Value = '_PMC_CROSSREF'

$ws.Range("E10").Value = '[Jie%Zhang%NULL%2,                    Huipeng%Lu%NULL%1,                    Haiping%Zeng%NULL%1,                    Shining%Zhang%NULL%1,                    Qifeng%Du%NULL%1,                    Tingyun%Jiang%NULL%1,                    Baoguo%Du%NULL%1]'
$ws.Range("I10").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E11").Value = '[Yu%Wang%NULL%1,                    Li-Yun%Zhu%NULL%2,                    Yu-Fen%Ma%NULL%1,                    Hai-Xin%Bo%NULL%2,                    Hai-Bo%Deng%NULL%1,                    Jing%Cao%NULL%2,                    Ying%Wang%NULL%0,                    Xiao-Jie%Wang%NULL%1,                    Yuan%Xu%NULL%1,                    Qiao-Dan%Lu%NULL%1,                    Hui%Wang%NULL%1,                    Xin-Juan%Wu%NULL%1]'
$ws.Range("I11").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E12").Value = '[Yanyu%Hu%NULL%1,                    Yingying%Chen%NULL%1,                    Yixiong%Zheng%NULL%1,                    Ciping%You%NULL%1,                    Jing%Tan%NULL%1,                    Lan%Hu%NULL%1,                    Zhenqing%Zhang%NULL%1,                    Lijun%Ding%NULL%1]'
$ws.Range("I12").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E13").Value = '[Yu-Fen%Ma%NULL%1,                    Wen%Li%NULL%0,                    Hai-Bao%Deng%NULL%1,                    Lei%Wang%NULL%0,                    Ying%Wang%NULL%0,                    Pei-Hong%Wang%NULL%1,                    Hai-Xin%Bo%NULL%0,                    Jing%Cao%NULL%0,                    Yu%Wang%NULL%0,                    Li-Yun%Zhu%NULL%0,                    Yuan%Yang%NULL%0,                    Teris%Cheung%NULL%0,                    Chee H.%Ng%NULL%1,                    Xinjuan%Wu%wuxinjuan@sina.com%0,                    Yu-Tao%Xiang%xyutly@gmail.com%0]'
$ws.Range("I13").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E14").Value = '[Jie%Zhang%NULL%0,                    Zhen%Yang%NULL%1,                    Xiao%Wang%NULL%0,                    Juan%Li%NULL%2,                    Lili%Dong%NULL%2,                    Lili%Dong%NULL%0,                    Fusheng%Wang%NULL%1,                    Yifei%Li%NULL%1,                    Ruihong%Wei%NULL%1,                    Jingping%Zhang%jpzhang1965@csu.edu.cn%1]'
$ws.Range("I14").Value = '_PMC_CROSSREF'

$ws.Range("E16").Value = '[Arman%Zarghami%xref no email%1,            Mojtaba%Farjam%xref no email%1,            Bahareh%Fakhraei%xref no email%1,            Kosar%Hashemzadeh%xref no email%1,            Mohammad Hosein%Yazdanpanah%xref no email%1]'
$ws.Range("I16").Value = '_CROSSREF'
